# [week8 day1] Work hard, study hard, and pay close attention to time to
# finish mathematics. The game has been deleted, but it feels a lot easier.
#
# Applies the content edits described by the commit:
#   - Row 37 (afternoon plan item): "学一下数学" -> "把单词背完"
#   - Row 38 (afternoon plan item): planned time window "15:10 - 14:00"
#     corrected to "15:10 - 16:00"
#   - Rows 36-38 actual-completion column (F) marked "yes"
#   - Selection / scroll position updated to reflect where the user was
#     working (A10 topmost, F41 selected)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Afternoon block (rows 36-38) content updates ---------------------
$ws.Range("C37").Value = "把单词背完"
$ws.Range("D38").Value = "15:10 - 16:00"

# Mark the three afternoon tasks as actually completed ("yes")
$ws.Range("F36").Value = "yes"
$ws.Range("F37").Value = "yes"
$ws.Range("F38").Value = "yes"

# --- View state: scroll/selection -------------------------------------
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("F41").Select()
